$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows after the header row (row 1), pushing existing data down.
$ws.Rows("2:7").Insert()

# Populate the newly inserted rows with the new test names (entered bottom-up,
# matching the original authoring order reflected in the shared-string table).
$ws.Range("A7").Value = "𝙁𝙍𝙀𝙀 𝘿𝙇_ Akon feat. Eminem - Smack That (Paranormila x SkaaR x E.B.A.H Edit)"
$ws.Range("A6").Value = "`$oho Bani - BLOCK THERAPIE (Techno Remix)"
$ws.Range("A5").Value = "IN FURCHT - THREEFOLD SYMMETRY"
$ws.Range("A4").Value = "DJ KINNƎR - CRAZY WIZARD SPELL"
$ws.Range("A3").Value = "TBK - KEVIN (TRANCESTRUDEL EDIT)"
$ws.Range("A2").Value = "RIANA HOLLEY & SERAFINA - KIM POSSIBLE"

# Update the selection to match the author's final cursor position.
$ws.Range("A13").Select()
